$d = $word.ActiveDocument

# --- 0. Drop the handful of unused styles that no longer ship with the
#        document (Heading2/4 and various leftover web-import styles).
#        Deleted from the highest original position down to the lowest so
#        that each style is still resolvable by name at the point its
#        Delete() is invoked. ---

$unusedStyles = @(
    "podcast-toolssubscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading4Char",
    "Heading2Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading4",
    "Heading2"
)
foreach ($styleName in $unusedStyles) {
    $style = $d.Styles.Item($styleName)
    $style.Delete()
}

# --- 1. First paragraph: add trailing spaces, then append a red-colored
#        parenthetical note in three runs (matches the target markup). ---

$p1 = $d.Paragraphs.Item(1)
$parEnd = $p1.Range.End - 1   # position just before the paragraph mark
$ip = $d.Range($parEnd, $parEnd)
$ip.InsertAfter("  ")

$p1 = $d.Paragraphs.Item(1)
$parEnd = $p1.Range.End - 1
$dash = [char]0x2013
$ip1 = $d.Range($parEnd, $parEnd)
$ip1.InsertAfter("(This is a change " + $dash + " Ve")
$ip1.Font.Color = 255

$p1 = $d.Paragraphs.Item(1)
$parEnd = $p1.Range.End - 1
$ip2 = $d.Range($parEnd, $parEnd)
$ip2.InsertAfter("rsion for main branch")
$ip2.Font.Color = 255

$p1 = $d.Paragraphs.Item(1)
$parEnd = $p1.Range.End - 1
$ip3 = $d.Range($parEnd, $parEnd)
$ip3.InsertAfter(")")
$ip3.Font.Color = 255

# --- 2. Remove the trailing "ank God almighty, we are free at last."
#        paragraph entirely, leaving "Shall be lifted-nevermore!" as the
#        last paragraph of the body. ---

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.Delete()

Write-Output "Paragraph 1: $($d.Paragraphs.Item(1).Range.Text)"
Write-Output "Last paragraph: $($d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)"
